# feat: add 2022-Q1 data
#
#  - Inserts a new "2022-Q1" worksheet (per-fund holding detail) right
#    before the "总计" (totals) sheet.
#  - Rebuilds the "总计" sheet with a new first data row summarizing
#    2022-Q1 (5 funds, 0.18 亿元) followed by the previously existing
#    quarters.
#
# NOTE: worksheet object variables in this COM host behave like
# positional handles rather than stable references: once the sheet
# collection is restructured (Copy/Add/Delete), a previously stored
# worksheet variable can end up pointing at whatever sheet now occupies
# its original slot. To stay safe, every worksheet reference below is
# re-fetched by index immediately before it is used.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 0. Remember the existing "总计" summary rows (date / count / value)
#    before rebuilding anything.
# ------------------------------------------------------------------
$total0 = $wb.Worksheets.Item($wb.Worksheets.Count)
$usedTotal = $total0.UsedRange
$totalRowCount = $usedTotal.Rows.Count

$existingDates = @()
$existingCounts = @()
$existingValues = @()
for ($r = 2; $r -le $totalRowCount; $r++) {
    $existingDates += , $total0.Cells.Item($r, 2).Value2
    $existingCounts += , $total0.Cells.Item($r, 3).Value2
    $existingValues += , $total0.Cells.Item($r, 4).Value2
}

# ------------------------------------------------------------------
# 1. Delete the original "总计" sheet. It currently holds the highest
#    sheetId, so removing it frees that id for reuse by the next
#    newly-created sheet (matching the target sheetId numbering).
# ------------------------------------------------------------------
$origTotal = $wb.Worksheets.Item($wb.Worksheets.Count)
$origTotal.Delete()

# ------------------------------------------------------------------
# 2. Build "2022-Q1" by copying the "2021-Q4" sheet (same 8-column
#    layout, styles and sheetPr/outline boilerplate), placed at the
#    end, then trimmed down to the required 6 rows (header + 5 funds).
# ------------------------------------------------------------------
$srcForQ1 = $wb.Worksheets.Item(3)
$endSheet1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcForQ1.Copy([System.Reflection.Missing]::Value, $endSheet1)

$q1 = $wb.Worksheets.Item(4)
$q1.Name = "2022-Q1"
$q1.Range("7:9").Delete()   # source sheet had 8 data rows, we only need 5

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$fundRows = @(
    @{ Code = "006836"; Name = "永赢惠泽一年定期开放灵活配置混合"; Scale = "8.57"; Position = "49.30"; Ratio = "0.63"; MarketValue = "0.0540"; Rank = 10 },
    @{ Code = "003594"; Name = "长盛盛崇灵活配置混合A";             Scale = "1.84"; Position = "43.63"; Ratio = "2.62"; MarketValue = "0.0482"; Rank = 5 },
    @{ Code = "080008"; Name = "长盛战略新兴产业灵活配置混合A";     Scale = "1.83"; Position = "40.84"; Ratio = "2.39"; MarketValue = "0.0437"; Rank = 5 },
    @{ Code = "001834"; Name = "长盛战略新兴产业灵活配置混合C";     Scale = "1.43"; Position = "40.84"; Ratio = "2.39"; MarketValue = "0.0342"; Rank = 5 },
    @{ Code = "003595"; Name = "长盛盛崇灵活配置混合C";             Scale = "0.11"; Position = "43.63"; Ratio = "2.62"; MarketValue = "0.0029"; Rank = 5 }
)

$r = 2
foreach ($fund in $fundRows) {
    $q1.Cells.Item($r, 1).Value = $r - 2

    # Columns B-G hold fund codes / decimal-looking strings that must stay
    # text (leading zeros, fixed decimal places). Force a temporary Text
    # format before assignment, then reset the style back to "Normal" so
    # no stray number-format style is left on the cell (matches the
    # un-styled inlineStr cells used by the sibling quarter sheets).
    $textRange = $q1.Range($q1.Cells.Item($r, 2), $q1.Cells.Item($r, 7))
    $textRange.NumberFormat = "@"
    $q1.Cells.Item($r, 2).Value = $fund.Code
    $q1.Cells.Item($r, 3).Value = $fund.Name
    $q1.Cells.Item($r, 4).Value = $fund.Scale
    $q1.Cells.Item($r, 5).Value = $fund.Position
    $q1.Cells.Item($r, 6).Value = $fund.Ratio
    $q1.Cells.Item($r, 7).Value = $fund.MarketValue
    $textRange.Style = "Normal"

    $q1.Cells.Item($r, 8).Value = $fund.Rank
    $r++
}

# ------------------------------------------------------------------
# 3. Rebuild "总计" the same way (copy from "2021-Q4" for matching
#    sheetPr/styles), trimmed to 4 columns x 5 rows, then repopulated
#    with the new 2022-Q1 row followed by the previous summary rows.
# ------------------------------------------------------------------
$srcForTotal = $wb.Worksheets.Item(3)
$endSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcForTotal.Copy([System.Reflection.Missing]::Value, $endSheet2)

$newTotal = $wb.Worksheets.Item(5)
$newTotal.Name = "总计"
$newTotal.Range("E:H").Delete()
$newTotal.Range("6:9").Delete()

$newTotal.Range("B1").Value = "日期"
$newTotal.Range("C1").Value = "持有数量(只)"
$newTotal.Range("D1").Value = "持有市值(亿元)"

$dates = @("2022-Q1") + $existingDates
$counts = @(5) + $existingCounts
$values = @(0.18) + $existingValues

for ($i = 0; $i -lt $dates.Length; $i++) {
    $rr = $i + 2
    $newTotal.Cells.Item($rr, 1).Value = $i

    $bRange = $newTotal.Cells.Item($rr, 2)
    $bRange.NumberFormat = "@"
    $newTotal.Cells.Item($rr, 2).Value = $dates[$i]
    $bRange.Style = "Normal"

    $newTotal.Cells.Item($rr, 3).Value = $counts[$i]
    $newTotal.Cells.Item($rr, 4).Value = $values[$i]
}
